$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) column F was refreshed for three
# rows that are duplicated across the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 5328
    $ws.Range("F6").Value = 167
    $ws.Range("F14").Value = 232
}
